$wb = $excel.ActiveWorkbook

# --- "About" sheet: document the EU adjustment ---------------------------
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("A23").Value = "To adjust this US study for the EU, we multiply by the ratio of EU:US LDVs (see file trans/BNVP for both the EU and US models)"

$wsAbout.Range("A24").Value = "2020 battery electric LDV, EU"
$wsAbout.Range("B24").Value = 54178.399285986576

$wsAbout.Range("A25").Value = "2020 battery electric LDV, EU"
$wsAbout.Range("B25").Value = 49995.669646960996

$wsAbout.Range("A24:A25").WrapText = $true
$wsAbout.Rows.Item(24).RowHeight = 58
$wsAbout.Rows.Item(25).RowHeight = 58

# --- "BRAaCTSC" sheet: adjust the US-sourced cost for the EU -------------
$wsBRA = $wb.Worksheets.Item("BRAaCTSC")
$wsBRA.Range("B2").Formula = "=Data!B6*(About!B24/About!B25)"

# --- Restore selections on each sheet (About stays the active tab) -------
[void]$wsBRA.Range("B3").Select()
[void]$wsAbout.Range("E33").Select()
